# Auto-generated edit script: "fixed monte carlo for ties"
$wb = $excel.ActiveWorkbook
$wsPlayoff = $wb.Worksheets.Item("Playoff Odds")
$wsRecord = $wb.Worksheets.Item("Record Odds")

# --- Record Odds sheet: update G (Expected_Final_Record) and H (Most_Likely_Record) text, and F (Playoff_Chance_Pct) ---
# Order matters: new shared strings get appended in first-use order, so we process
# all of column G (rows 2-11) first, then all of column H (rows 2-11), to reproduce
# the canonical shared string table ordering.
$wsRecord.Cells.Item(2, 7).Value = "9.3-4.7-0.1"
$wsRecord.Cells.Item(3, 7).Value = "8.5-5.4-0.1"
$wsRecord.Cells.Item(4, 7).Value = "7.9-6.1-0.1"
$wsRecord.Cells.Item(5, 7).Value = "7.2-6.8"
$wsRecord.Cells.Item(6, 7).Value = "7.0-7.0-0.1"
$wsRecord.Cells.Item(7, 7).Value = "6.8-7.1-0.1"
$wsRecord.Cells.Item(8, 7).Value = "6.4-7.5-0.1"
$wsRecord.Cells.Item(9, 7).Value = "6.5-7.4-0.1"
$wsRecord.Cells.Item(10, 7).Value = "5.3-8.6-0.1"
$wsRecord.Cells.Item(11, 7).Value = "4.8-9.2-0.1"
$wsRecord.Cells.Item(2, 8).Value = "10-4"
$wsRecord.Cells.Item(3, 8).Value = "9-5"
$wsRecord.Cells.Item(4, 8).Value = "8-6"
$wsRecord.Cells.Item(5, 8).Value = "7-7"
$wsRecord.Cells.Item(6, 8).Value = "7-7"
$wsRecord.Cells.Item(7, 8).Value = "7-7"
$wsRecord.Cells.Item(8, 8).Value = "7-7"
$wsRecord.Cells.Item(9, 8).Value = "7-7"
$wsRecord.Cells.Item(10, 8).Value = "5-9"
$wsRecord.Cells.Item(11, 8).Value = "4-10"

# --- Record Odds sheet: update F (Playoff_Chance_Pct) ---
$wsRecord.Cells.Item(2, 6).Value = 97.39999999999999
$wsRecord.Cells.Item(3, 6).Value = 90.3
$wsRecord.Cells.Item(4, 6).Value = 80.2
$wsRecord.Cells.Item(5, 6).Value = 71.8
$wsRecord.Cells.Item(6, 6).Value = 63
$wsRecord.Cells.Item(7, 6).Value = 57.59999999999999
$wsRecord.Cells.Item(8, 6).Value = 55.60000000000001
$wsRecord.Cells.Item(9, 6).Value = 50.3
$wsRecord.Cells.Item(10, 6).Value = 22.4
$wsRecord.Cells.Item(11, 6).Value = 11.4

# --- Playoff Odds sheet: update the full percentage grid (B:L, rows 2-11) ---
$wsPlayoff.Cells.Item(2, 2).Value = 47.2
$wsPlayoff.Cells.Item(2, 3).Value = 23.1
$wsPlayoff.Cells.Item(2, 4).Value = 11.9
$wsPlayoff.Cells.Item(2, 5).Value = 8.1
$wsPlayoff.Cells.Item(2, 6).Value = 3.7
$wsPlayoff.Cells.Item(2, 7).Value = 3.4
$wsPlayoff.Cells.Item(2, 8).Value = 1.3
$wsPlayoff.Cells.Item(2, 9).Value = 0.9
$wsPlayoff.Cells.Item(2, 11).Value = 0.1
$wsPlayoff.Cells.Item(2, 12).Value = 97.4
$wsPlayoff.Cells.Item(3, 2).Value = 23.3
$wsPlayoff.Cells.Item(3, 3).Value = 22.3
$wsPlayoff.Cells.Item(3, 4).Value = 16.1
$wsPlayoff.Cells.Item(3, 5).Value = 13.4
$wsPlayoff.Cells.Item(3, 6).Value = 8.5
$wsPlayoff.Cells.Item(3, 7).Value = 6.7
$wsPlayoff.Cells.Item(3, 8).Value = 4.4
$wsPlayoff.Cells.Item(3, 9).Value = 2.7
$wsPlayoff.Cells.Item(3, 10).Value = 2
$wsPlayoff.Cells.Item(3, 11).Value = 0.6
$wsPlayoff.Cells.Item(3, 12).Value = 90.3
$wsPlayoff.Cells.Item(4, 2).Value = 8.6
$wsPlayoff.Cells.Item(4, 3).Value = 15.8
$wsPlayoff.Cells.Item(4, 4).Value = 18.2
$wsPlayoff.Cells.Item(4, 5).Value = 13.2
$wsPlayoff.Cells.Item(4, 6).Value = 13.2
$wsPlayoff.Cells.Item(4, 7).Value = 11.2
$wsPlayoff.Cells.Item(4, 8).Value = 7.9
$wsPlayoff.Cells.Item(4, 9).Value = 7
$wsPlayoff.Cells.Item(4, 10).Value = 4.1
$wsPlayoff.Cells.Item(4, 11).Value = 0.8
$wsPlayoff.Cells.Item(4, 12).Value = 80.2
$wsPlayoff.Cells.Item(5, 2).Value = 7.2
$wsPlayoff.Cells.Item(5, 3).Value = 11.4
$wsPlayoff.Cells.Item(5, 5).Value = 12.2
$wsPlayoff.Cells.Item(5, 6).Value = 14.4
$wsPlayoff.Cells.Item(5, 7).Value = 13.3
$wsPlayoff.Cells.Item(5, 8).Value = 10.2
$wsPlayoff.Cells.Item(5, 9).Value = 8.3
$wsPlayoff.Cells.Item(5, 10).Value = 6
$wsPlayoff.Cells.Item(5, 11).Value = 3.7
$wsPlayoff.Cells.Item(5, 12).Value = 71.8
$wsPlayoff.Cells.Item(6, 2).Value = 4.1
$wsPlayoff.Cells.Item(6, 3).Value = 8.9
$wsPlayoff.Cells.Item(6, 4).Value = 12.3
$wsPlayoff.Cells.Item(6, 5).Value = 11.8
$wsPlayoff.Cells.Item(6, 6).Value = 11.8
$wsPlayoff.Cells.Item(6, 7).Value = 14.1
$wsPlayoff.Cells.Item(6, 8).Value = 12.5
$wsPlayoff.Cells.Item(6, 9).Value = 11.5
$wsPlayoff.Cells.Item(6, 10).Value = 8.4
$wsPlayoff.Cells.Item(6, 11).Value = 4.6
$wsPlayoff.Cells.Item(6, 12).Value = 63
$wsPlayoff.Cells.Item(7, 2).Value = 3.8
$wsPlayoff.Cells.Item(7, 3).Value = 6.4
$wsPlayoff.Cells.Item(7, 4).Value = 9.4
$wsPlayoff.Cells.Item(7, 5).Value = 13.6
$wsPlayoff.Cells.Item(7, 6).Value = 11.8
$wsPlayoff.Cells.Item(7, 7).Value = 12.6
$wsPlayoff.Cells.Item(7, 8).Value = 13.4
$wsPlayoff.Cells.Item(7, 9).Value = 13
$wsPlayoff.Cells.Item(7, 10).Value = 10.2
$wsPlayoff.Cells.Item(7, 11).Value = 5.8
$wsPlayoff.Cells.Item(7, 12).Value = 57.6
$wsPlayoff.Cells.Item(8, 2).Value = 3.9
$wsPlayoff.Cells.Item(8, 3).Value = 6.9
$wsPlayoff.Cells.Item(8, 4).Value = 8.8
$wsPlayoff.Cells.Item(8, 5).Value = 11.6
$wsPlayoff.Cells.Item(8, 6).Value = 13.1
$wsPlayoff.Cells.Item(8, 7).Value = 11.3
$wsPlayoff.Cells.Item(8, 8).Value = 14.4
$wsPlayoff.Cells.Item(8, 9).Value = 11.2
$wsPlayoff.Cells.Item(8, 10).Value = 11.2
$wsPlayoff.Cells.Item(8, 11).Value = 7.6
$wsPlayoff.Cells.Item(8, 12).Value = 55.6
$wsPlayoff.Cells.Item(9, 2).Value = 1.8
$wsPlayoff.Cells.Item(9, 3).Value = 3.8
$wsPlayoff.Cells.Item(9, 4).Value = 7.1
$wsPlayoff.Cells.Item(9, 5).Value = 10
$wsPlayoff.Cells.Item(9, 6).Value = 13.6
$wsPlayoff.Cells.Item(9, 7).Value = 14
$wsPlayoff.Cells.Item(9, 8).Value = 16.5
$wsPlayoff.Cells.Item(9, 9).Value = 15.1
$wsPlayoff.Cells.Item(9, 10).Value = 11.6
$wsPlayoff.Cells.Item(9, 11).Value = 6.5
$wsPlayoff.Cells.Item(9, 12).Value = 50.3
$wsPlayoff.Cells.Item(10, 2).Value = 0
$wsPlayoff.Cells.Item(10, 3).Value = 0.9
$wsPlayoff.Cells.Item(10, 4).Value = 1.8
$wsPlayoff.Cells.Item(10, 5).Value = 4.1
$wsPlayoff.Cells.Item(10, 6).Value = 6.2
$wsPlayoff.Cells.Item(10, 7).Value = 9.4
$wsPlayoff.Cells.Item(10, 8).Value = 11.2
$wsPlayoff.Cells.Item(10, 9).Value = 17.4
$wsPlayoff.Cells.Item(10, 10).Value = 25.1
$wsPlayoff.Cells.Item(10, 11).Value = 23.9
$wsPlayoff.Cells.Item(10, 12).Value = 22.4
$wsPlayoff.Cells.Item(11, 2).Value = 0.1
$wsPlayoff.Cells.Item(11, 3).Value = 0.5
$wsPlayoff.Cells.Item(11, 4).Value = 1.1
$wsPlayoff.Cells.Item(11, 6).Value = 3.7
$wsPlayoff.Cells.Item(11, 7).Value = 4
$wsPlayoff.Cells.Item(11, 8).Value = 8.2
$wsPlayoff.Cells.Item(11, 9).Value = 12.9
$wsPlayoff.Cells.Item(11, 10).Value = 21.1
$wsPlayoff.Cells.Item(11, 11).Value = 46.4
$wsPlayoff.Cells.Item(11, 12).Value = 11.4

Write-Host "done"
